$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - match formatting of the other header cells (bold + border, style index 1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data cells F2:F25 - time_taken values (kept as plain text)
$ws.Range("F2").Value = "2021-10-05 10:51:49.178008"
$ws.Range("F3").Value = "2021-10-05 10:51:49.178019"
$ws.Range("F4").Value = "2021-10-05 10:51:49.178023"
$ws.Range("F5").Value = "2021-10-05 10:51:49.178025"
$ws.Range("F6").Value = "2021-10-05 10:51:49.178028"
$ws.Range("F7").Value = "2021-10-05 10:51:49.178031"
$ws.Range("F8").Value = "2021-10-05 10:51:49.178034"
$ws.Range("F9").Value = "2021-10-05 10:51:49.178036"
$ws.Range("F10").Value = "2021-10-05 10:51:49.178039"
$ws.Range("F11").Value = "2021-10-05 10:51:49.178042"
$ws.Range("F12").Value = "2021-10-05 10:51:49.178044"
$ws.Range("F13").Value = "2021-10-05 10:51:49.178047"
$ws.Range("F14").Value = "2021-10-05 10:51:49.178049"
$ws.Range("F15").Value = "2021-10-05 10:51:49.178052"
$ws.Range("F16").Value = "2021-10-05 10:51:49.178054"
$ws.Range("F17").Value = "2021-10-05 10:51:49.178057"
$ws.Range("F18").Value = "2021-10-05 10:51:49.178060"
$ws.Range("F19").Value = "2021-10-05 10:51:49.178062"
$ws.Range("F20").Value = "2021-10-05 10:51:49.178065"
$ws.Range("F21").Value = "2021-10-05 10:51:49.178068"
$ws.Range("F22").Value = "2021-10-05 10:51:49.178070"
$ws.Range("F23").Value = "2021-10-05 10:51:49.178073"
$ws.Range("F24").Value = "2021-10-05 10:51:49.178075"
$ws.Range("F25").Value = "2021-10-05 10:51:49.178078"
